$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values in this sheet are stored as plain text (not numbers), even when
# they look numeric (e.g. "4.51"), because the source values use a dot as a
# thousands/price separator rather than a decimal point, and trailing zeros
# must be preserved exactly as authored. For any replacement text that Excel
# would otherwise auto-detect as a number, force a text quote-prefix first.

$ws.Range("D2").Value = '68.310.73'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '2.715.66'
$ws.Range("E3").Value = '  +2.61%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''607.08'
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("D6").Value = '''166.90'
$ws.Range("E6").Value = '  +4.97%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.555'
$ws.Range("E8").Value = '  +3.07%  '
$ws.Range("D9").Value = '2.714.43'
$ws.Range("E9").Value = '  +2.60%  '
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '''0.158'
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '''0.365'
$ws.Range("E12").Value = '  +4.18%  '
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").Value = '''28.51'
$ws.Range("E14").Value = '  +2.13%  '
$ws.Range("D15").Value = '3.217.82'
$ws.Range("E15").Value = '  +2.76%  '
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '68.264.31'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '2.722.28'
$ws.Range("E18").Value = '  +3.16%  '
$ws.Range("D19").Value = '''11.87'
$ws.Range("E19").Value = '  +4.01%  '
$ws.Range("D20").Value = '''371.55'
$ws.Range("E20").Value = '  +2.46%  '
$ws.Range("D21").Value = '''7.64'
$ws.Range("E21").Value = '  +3.38%  '
$ws.Range("D22").Value = '''4.51'
$ws.Range("E22").Value = '  +2.44%  '
$ws.Range("D23").Value = '''4.97'
$ws.Range("E23").Value = '  +4.30%  '
$ws.Range("D24").Value = '''2.08'
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("D25").Value = '''73.02'
$ws.Range("E25").Value = '  -2.00%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '''10.19'
$ws.Range("E27").Value = '  +4.12%  '
$ws.Range("D28").Value = '2.854.69'
$ws.Range("E28").Value = '  +2.75%  '
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '''578.54'
$ws.Range("E31").Value = '  +2.86%  '
$ws.Range("D32").Value = '''8.20'
$ws.Range("E32").Value = '  +1.89%  '
$ws.Range("E33").Value = '  +2.27%  '
$ws.Range("E34").Value = '  +6.96%  '
$ws.Range("E35").Value = '  +2.05%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''1.60'
$ws.Range("E37").Value = '  -3.62%  '
$ws.Range("D38").Value = '''162.13'
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("D39").Value = '''19.88'
$ws.Range("E39").Value = '  +1.41%  '
$ws.Range("E40").Value = '  +2.03%  '
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("D42").Value = '''5.41'
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("D43").Value = '''17.99'
$ws.Range("E43").Value = '  +0.98%  '
$ws.Range("D44").Value = '''2.62'
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = '0.0₆0310'
$ws.Range("E46").Value = '  -3.04%  '
$ws.Range("E47").Value = '  +1.48%  '
$ws.Range("D48").Value = '''0.595'
$ws.Range("E48").Value = '  +4.00%  '
$ws.Range("D49").Value = '''154.85'
$ws.Range("E49").Value = '  -1.92%  '
$ws.Range("E50").Value = '  +2.52%  '
$ws.Range("D51").Value = '''1.77'
$ws.Range("E51").Value = '  +5.01%  '
